$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.207.93'
$ws.Range("E2").Value = '  -4.33%  '
$ws.Range("D3").Value = '2.269.05'
$ws.Range("E3").Value = '  -6.03%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.36'
$ws.Range("E5").Value = '  -4.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '80.98'
$ws.Range("E6").Value = '  -8.60%  '
$ws.Range("E7").Value = '  -4.79%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.463'
$ws.Range("E9").Value = '  -6.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0778'
$ws.Range("E10").Value = '  -6.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.22'
$ws.Range("E11").Value = '  -10.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.05'
$ws.Range("E12").Value = '  -11.81%  '
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").Value = '2.625.34'
$ws.Range("E14").Value = '  -5.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.11'
$ws.Range("E15").Value = '  -10.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.21'
$ws.Range("E16").Value = '  -9.16%  '
$ws.Range("D17").Value = '2.278.45'
$ws.Range("E17").Value = '  -5.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.724'
$ws.Range("E18").Value = '  -5.94%  '
$ws.Range("D19").Value = '39.141.29'
$ws.Range("E19").Value = '  -4.07%  '
$ws.Range("D20").Value = '0.0₃0867'
$ws.Range("E20").Value = '  -5.76%  '
$ws.Range("E21").Value = '  -6.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.29'
$ws.Range("E22").Value = '  -6.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.04'
$ws.Range("E23").Value = '  -7.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '228.55'
$ws.Range("E24").Value = '  -2.13%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  -9.54%  '
$ws.Range("E27").Value = '  -5.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.46'
$ws.Range("E28").Value = '  -6.37%  '
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.90'
$ws.Range("E30").Value = '  -6.52%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.99'
$ws.Range("E31").Value = '  -5.43%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.77'
$ws.Range("E32").Value = '  -6.94%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.84'
$ws.Range("E34").Value = '  -8.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.37'
$ws.Range("E35").Value = '  -3.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0688'
$ws.Range("E36").Value = '  -6.65%  '
$ws.Range("E37").Value = '  -3.26%  '
$ws.Range("E38").Value = '  -7.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0955'
$ws.Range("E39").Value = '  -4.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.85'
$ws.Range("E40").Value = '  -9.22%  '
$ws.Range("E41").Value = '  -8.38%  '
$ws.Range("E42").Value = '  -6.02%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.938.13'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.23'
$ws.Range("E44").Value = '  -3.03%  '
$ws.Range("E45").Value = '  -7.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.62'
$ws.Range("E46").Value = '  -8.97%  '
$ws.Range("E47").Value = '  -3.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.54'
$ws.Range("E48").Value = '  -10.76%  '
$ws.Range("D49").Value = '2.498.81'
$ws.Range("E49").Value = '  -5.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.78'
$ws.Range("E50").Value = '  -5.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '66.77'
$ws.Range("E51").Value = '  -8.79%  '

# Reset number format on cells we forced to text, so no residual style marker remains
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
